# "Added intro to course slides."
# Insert two new columns ("Youtube" and "Slides") in front of the existing
# "Resources" column on the Meetups sheet, and record the slide deck for
# the first ("Intro to the Course") session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Insert two blank columns before column F ("Resources"); this pushes the
# old F column (and its cells/formatting) to H.
$ws.Columns("F:G").Insert()

# Header row for the two new columns.
$ws.Range("F1").Value = "Youtube"
$ws.Range("G1").Value = "Slides"

# Slide deck reference for the "Intro to the Course" meetup (row 2).
$ws.Range("G2").Value = "00-Intro_to_Course"

# Size the new columns.
$ws.Columns("G").ColumnWidth = 16.1666666666667
$ws.Columns("H").ColumnWidth = 22.3307291666667

# Leave the selection on the newly added Youtube column header.
$ws.Range("F2").Select()
